$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark FILL instruction (and neighbouring rows) as supported in the AST /
# Parser columns of the checklist.

# LDM (row 27) and STM (row 28): mark AST column (D) as done
$ws.Range("D27").Value = "y"
$ws.Range("D28").Value = "y"

# DCD (row 31), EQU (row 32), FILL (row 33), END (row 34):
# mark Parser (C) and AST (D) columns as done
$ws.Range("C31").Value = "y"
$ws.Range("D31").Value = "y"
$ws.Range("C32").Value = "y"
$ws.Range("D32").Value = "y"
$ws.Range("C33").Value = "y"
$ws.Range("D33").Value = "y"
$ws.Range("C34").Value = "y"
$ws.Range("D34").Value = "y"

# Update the active cell selection to reflect where the edit was made
$ws.Range("D28").Select()
